# This script re-shuffles the course-syllabus content of LOB1272.docx so that:
#  - the "Objetivos" section now holds the (PT/EN) description of water/sewage
#    system planning & operation (previously under "Programa resumido"),
#  - the "Docente(s) Responsavel(eis)" bullet gets the text that used to be the
#    "Objetivos" (PT) paragraph,
#  - "Programa resumido" gets the detailed "Programa" (PT) text, and its EN
#    (italic) companion becomes the old "Objetivos" EN sentence,
#  - "Programa" gets the short evaluation-method sentence that used to open
#    "Avaliacao",
#  - inside "Avaliacao", the Metodo / Criterio / Norma de recuperacao values
#    each shift to the next slot, with "Norma de recuperacao" ending up with
#    the whole bibliography text block,
#  - the "Bibliografia" heading moves down (now right above the final bullet),
#    and the paragraph that used to hold the bibliography block becomes the
#    professor bullet ("7455355 - Robson da Silva Rocha").
#
# Paragraph styles / run formatting (bold "label" runs, italic EN blurbs,
# ListBullet styles, line breaks) all stay exactly where they are; only the
# w:t text content moves between slots. To make this robust regardless of
# execution order we:
#   1. Read every piece of text that needs to move into a variable FIRST,
#      before any mutation happens (so reads can never observe a half-done
#      edit).
#   2. Overwrite the three "Avaliacao" run slots using absolute character
#      offsets captured from the pristine document, processing them from the
#      end of the paragraph backwards so writing a (possibly different
#      length) replacement never invalidates an offset we still need.
#   3. Overwrite the remaining (single-run) paragraphs by indexed
#      Paragraphs.Item(n).Range, which stays valid no matter how earlier
#      paragraphs' text lengths changed.

$d = $word.ActiveDocument

function Get-MatchText($searchText) {
    $dup = $d.Content.Duplicate
    $ok = $dup.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $searchText"
    }
    return $dup.Text
}

function Get-MatchRange($searchText) {
    $dup = $d.Content.Duplicate
    $ok = $dup.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $searchText"
    }
    return $dup
}

# ---- 1. Read all source values (before any mutation) ----------------------

$A = Get-MatchText("Apresentar aos alunos uma visão integrada dos principais componentes dos sistemas de abastecimento de água e coleta esgotos sanitários")
$B = Get-MatchText("Provide to the students an integrated view on key components of water supply and wastewater systems")
$C = Get-MatchText("7455355 - Robson da Silva Rocha")
$D = Get-MatchText("Planejamento, projeto e operação de sistemas de abastecimento de água e de esgotos: captação, adução, reservação e distribuição, redes coletoras, elevatórias, emissários e interceptores")
$E = Get-MatchText("Planning, design and operation of water supply and wastewater systems: sources, water mains, reservation and distribution, wastewater collection, pumping stations, interceptors and emissaries.")
$F = Get-MatchText("Sistemas de abastecimento de água. Partes constitutivas do sistema. Captações de águas superficiais e subterrâneas. Estações elevatórias e adutoras. Reservatórios. Redes de distribuição de água. Controle de perdas. Sistemas de esgotos sanitários. Rede coletora. Interceptores e emissários. Sifão Invertido. Estações elevatórias e linhas de recalque.")

# Character ranges (+ text) of the three "value" runs inside the Avaliacao
# bullet (captured now, while the document is still pristine).
$rG = Get-MatchRange("Avaliação baseada em provas, exercícios, trabalhos práticos e relatórios")
$gStart = $rG.Start
$gEnd = $rG.End
$G = $rG.Text

$rH = Get-MatchRange("Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios")
$hStart = $rH.Start
$hEnd = $rH.End
$H = $rH.Text

$rI = Get-MatchRange("1 (uma) prova de recuperação (R), sendo considerado aprovado se 0,5(NF + R) >= 5,0")
$iStart = $rI.Start
$iEnd = $rI.End
$I = $rI.Text

# Whole bibliography paragraph (multiple runs joined by soft line breaks);
# Range.Text includes the trailing paragraph mark (Chr 13) which we drop.
$p19 = $d.Paragraphs.Item(19)
$J = $p19.Range.Text
if ($J.Length -gt 0 -and [int][char]$J.Substring($J.Length - 1, 1) -eq 13) {
    $J = $J.Substring(0, $J.Length - 1)
}

# ---- 2. Overwrite the Avaliacao run slots, back-to-front ------------------

$d.Range($iStart, $iEnd).Text = $J
$d.Range($hStart, $hEnd).Text = $I
$d.Range($gStart, $gEnd).Text = $H

# ---- 3. Overwrite the remaining single-run paragraphs by index ------------

$d.Paragraphs.Item(6).Range.Text = $D
$d.Paragraphs.Item(7).Range.Text = $E
$d.Paragraphs.Item(9).Range.Text = $A
$d.Paragraphs.Item(11).Range.Text = $F
$d.Paragraphs.Item(12).Range.Text = $B
$d.Paragraphs.Item(14).Range.Text = $G
$d.Paragraphs.Item(19).Range.Text = $C
